$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.628473391421153
$ws.Range("C2").Value = 0.5925921763866313
$ws.Range("D2").Value = 0.6605252660765188
$ws.Range("E2").Value = 0.2694002117863263
$ws.Range("G2").Value = 1.618817807492434
$ws.Range("H2").Value = 1.3914997392651
$ws.Range("I2").Value = 0.9754574343271756
$ws.Range("J2").Value = 0.1402073362942957
$ws.Range("N2").Value = 1.314684938046567
$ws.Range("B3").Value = 1.504356375888165
$ws.Range("C3").Value = 0.5463184728849342
$ws.Range("D3").Value = 0.6507702766687657
$ws.Range("E3").Value = 0.2645580228302933
$ws.Range("G3").Value = 1.593077908450823
$ws.Range("H3").Value = 1.386880558471887
$ws.Range("I3").Value = 0.9749819146558849
$ws.Range("J3").Value = 0.1369585051028821
$ws.Range("N3").Value = 1.334350938428487
$ws.Range("B4").Value = 1.428850933360536
$ws.Range("C4").Value = 0.5181917515143937
$ws.Range("D4").Value = 0.6451483806956162
$ws.Range("E4").Value = 0.2617418316259048
$ws.Range("G4").Value = 1.578514329356182
$ws.Range("H4").Value = 1.384899098715238
$ws.Range("I4").Value = 0.9753654369932434
$ws.Range("J4").Value = 0.1350503228013196
$ws.Range("N4").Value = 1.347009504765835
$ws.Range("B5").Value = 1.398257842371095
$ws.Range("C5").Value = 0.5068009986992479
$ws.Range("D5").Value = 0.6429496166272202
$ws.Range("E5").Value = 0.2606335290419679
$ws.Range("G5").Value = 1.572889633524113
$ws.Range("H5").Value = 1.384305833629156
$ws.Range("I5").Value = 0.9756908265652768
$ws.Range("J5").Value = 0.134294394602442
$ws.Range("N5").Value = 1.352314462717342
$ws.Range("B6").Value = 1.393188490800583
$ws.Range("C6").Value = 0.504913848392448
$ws.Range("D6").Value = 0.6425900750007543
$ws.Range("E6").Value = 0.2604518668217537
$ws.Range("G6").Value = 1.571974328982606
$ws.Range("H6").Value = 1.384220238026984
$ws.Range("I6").Value = 0.9757550485303526
$ws.Range("J6").Value = 0.1341701790890042
$ws.Range("N6").Value = 1.353204187278902
$ws.Range("B7").Value = 1.428437632653356
$ws.Range("C7").Value = 0.5180378447853968
$ws.Range("D7").Value = 0.6451183543991306
$ws.Range("E7").Value = 0.2617267256586473
$ws.Range("G7").Value = 1.578437219798644
$ws.Range("H7").Value = 1.384890231451749
$ws.Range("I7").Value = 0.9753691416388719
$ws.Range("J7").Value = 0.1350400404824512
$ws.Range("N7").Value = 1.347080456500355
$ws.Range("B8").Value = 1.585531447687913
$ws.Range("C8").Value = 0.5765773148849576
$ws.Range("D8").Value = 0.6570852661189406
$ws.Range("E8").Value = 0.2676979642985842
$ws.Range("G8").Value = 1.609683976118077
$ws.Range("H8").Value = 1.389729195255029
$ws.Range("I8").Value = 0.9751528158736349
$ws.Range("J8").Value = 0.1390690959313972
$ws.Range("N8").Value = 1.321344505696145
$ws.Range("B9").Value = 1.899222039286201
$ws.Range("C9").Value = 0.693676002336133
$ws.Range("D9").Value = 0.6834838731692514
$ws.Range("E9").Value = 0.2806602570430812
$ws.Range("G9").Value = 1.680897328382883
$ws.Range("H9").Value = 1.406036198120773
$ws.Range("I9").Value = 0.9801240946056708
$ws.Range("J9").Value = 0.1476628058938587
$ws.Range("N9").Value = 1.275517280668474
$ws.Range("B10").Value = 2.133222665511539
$ws.Range("C10").Value = 0.7811740183080929
$ws.Range("D10").Value = 0.7046886589652388
$ws.Range("E10").Value = 0.2909592655093363
$ws.Range("G10").Value = 1.739414207452
$ws.Range("H10").Value = 1.422226896527576
$ws.Range("I10").Value = 0.9871175613414493
$ws.Range("J10").Value = 0.1544074671669904
$ws.Range("N10").Value = 1.244690768785377
$ws.Range("B11").Value = 2.240466272806088
$ws.Range("C11").Value = 0.82131148431273
$ws.Range("D11").Value = 0.7147334137535779
$ws.Range("E11").Value = 0.2958156818574977
$ws.Range("G11").Value = 1.767411088881744
$ws.Range("H11").Value = 1.430518622586789
$ws.Range("I11").Value = 0.9910360179368638
$ws.Range("J11").Value = 0.1575712757594658
$ws.Range("N11").Value = 1.231287276134942
$ws.Range("B12").Value = 2.28119246360427
$ws.Range("C12").Value = 0.8365595230486065
$ws.Range("D12").Value = 0.718594789552327
$ws.Range("E12").Value = 0.2976795171883495
$ws.Range("G12").Value = 1.778213190359565
$ws.Range("H12").Value = 1.433792645245006
$ws.Range("I12").Value = 0.9926267418299943
$ws.Range("J12").Value = 0.1587832215649883
$ws.Range("N12").Value = 1.22630113180536
$ws.Range("B13").Value = 2.272416203065745
$ws.Range("C13").Value = 0.8332733995184753
$ws.Range("D13").Value = 0.7177606040892215
$ws.Range("E13").Value = 0.2972770003577949
$ws.Range("G13").Value = 1.775877820562698
$ws.Range("H13").Value = 1.433081546149083
$ws.Range("I13").Value = 0.9922793826429199
$ws.Range("J13").Value = 0.1585215879490676
$ws.Range("N13").Value = 1.227371002218768
$ws.Range("B14").Value = 2.24381452328879
$ws.Range("C14").Value = 0.8225649659490841
$ws.Range("D14").Value = 0.7150499342609749
$ws.Range("E14").Value = 0.295968522339912
$ws.Range("G14").Value = 1.768295757215185
$ws.Range("H14").Value = 1.430785285127598
$ws.Range("I14").Value = 0.9911647403515076
$ws.Range("J14").Value = 0.1576707043436016
$ws.Range("N14").Value = 1.230875267864832
$ws.Range("B15").Value = 2.226310221741471
$ws.Range("C15").Value = 0.8160121270581726
$ws.Range("D15").Value = 0.713397087703953
$ws.Range("E15").Value = 0.2951702788355703
$ws.Range("G15").Value = 1.763677675201222
$ws.Range("H15").Value = 1.429396252979558
$ws.Range("I15").Value = 0.9904959367169894
$ws.Range("J15").Value = 0.1571513250408145
$ws.Range("N15").Value = 1.233033394005462
$ws.Range("B16").Value = 2.126230133618265
$ws.Range("C16").Value = 0.7785577327544502
$ws.Range("D16").Value = 0.7040402591664758
$ws.Range("E16").Value = 0.2906453501989716
$ws.Range("G16").Value = 1.737612450368971
$ws.Range("H16").Value = 1.421703732531057
$ws.Range("I16").Value = 0.9868763895463104
$ws.Range("J16").Value = 0.1542026400863392
$ws.Range("N16").Value = 1.245579218212794
$ws.Range("B17").Value = 2.065038674196501
$ws.Range("C17").Value = 0.7556668275727247
$ws.Range("D17").Value = 0.6984024677050797
$ws.Range("E17").Value = 0.2879134632094846
$ws.Range("G17").Value = 1.721976614360784
$ws.Range("H17").Value = 1.417222551740167
$ws.Range("I17").Value = 0.9848453359952885
$ws.Range("J17").Value = 0.1524183010471347
$ws.Range("N17").Value = 1.25343463605319
$ws.Range("B18").Value = 2.029917785779617
$ws.Range("C18").Value = 0.742531990154589
$ws.Range("D18").Value = 0.6951972563643096
$ws.Range("E18").Value = 0.286358271154306
$ws.Range("G18").Value = 1.713112793468355
$ws.Range("H18").Value = 1.414732246762554
$ws.Range("I18").Value = 0.9837464521746142
$ws.Range("J18").Value = 0.1514009895638111
$ws.Range("N18").Value = 1.25801116599125
$ws.Range("B19").Value = 2.018039278181789
$ws.Range("C19").Value = 0.7380901298934077
$ws.Range("D19").Value = 0.6941184571871872
$ws.Range("E19").Value = 0.2858344726539173
$ws.Range("G19").Value = 1.710133829645457
$ws.Range("H19").Value = 1.413904013571198
$ws.Range("I19").Value = 0.9833862684046082
$ws.Range("J19").Value = 0.1510580860670672
$ws.Range("N19").Value = 1.259570702766053
$ws.Range("B20").Value = 2.071544859005144
$ws.Range("C20").Value = 0.7581003474140857
$ws.Range("D20").Value = 0.6989987370621407
$ws.Range("E20").Value = 0.2882026079197786
$ws.Range("G20").Value = 1.723627657057534
$ws.Range("H20").Value = 1.417690554769734
$ws.Range("I20").Value = 0.9850543638409519
$ws.Range("J20").Value = 0.152607315302518
$ws.Range("N20").Value = 1.252592376841394
$ws.Range("B21").Value = 2.252212391105616
$ws.Range("C21").Value = 0.8257089602074643
$ws.Range("D21").Value = 0.7158445565671627
$ws.Range("E21").Value = 0.2963521791305581
$ws.Range("G21").Value = 1.770517339184124
$ws.Range("H21").Value = 1.431456105287651
$ws.Range("I21").Value = 0.9914892296624203
$ws.Range("J21").Value = 0.1579202517410891
$ws.Range("N21").Value = 1.229843549034509
$ws.Range("B22").Value = 2.37096197558202
$ws.Range("C22").Value = 0.8701802237956713
$ws.Range("D22").Value = 0.7271904317157123
$ws.Range("E22").Value = 0.3018231255369415
$ws.Range("G22").Value = 1.802330962337066
$ws.Range("H22").Value = 1.441234904007842
$ws.Range("I22").Value = 0.9963182440432661
$ws.Range("J22").Value = 0.1614735285947404
$ws.Range("N22").Value = 1.215497535000684
$ws.Range("B23").Value = 2.307521186695794
$ws.Range("C23").Value = 0.8464187162177268
$ws.Range("D23").Value = 0.7211040555385182
$ws.Range("E23").Value = 0.2988898762844272
$ws.Range("G23").Value = 1.785243766233435
$ws.Range("H23").Value = 1.435943902006812
$ws.Range("I23").Value = 0.9936835625993155
$ws.Range("J23").Value = 0.1595696262561574
$ws.Range("N23").Value = 1.223106418200292
$ws.Range("B24").Value = 2.068603229475116
$ws.Range("C24").Value = 0.7570000738858198
$ws.Range("D24").Value = 0.6987290514845483
$ws.Range("E24").Value = 0.2880718376228231
$ws.Range("G24").Value = 1.722880830138337
$ws.Range("H24").Value = 1.417478702842232
$ws.Range("I24").Value = 0.9849596481355931
$ws.Range("J24").Value = 0.1525218354020836
$ws.Range("N24").Value = 1.252972974295447
$ws.Range("B25").Value = 1.813747498050532
$ws.Range("C25").Value = 0.6617448455531303
$ws.Range("D25").Value = 0.6760262026383259
$ws.Range("E25").Value = 0.2770182822583962
$ws.Range("G25").Value = 1.660554111595957
$ws.Range("H25").Value = 1.400889336373922
$ws.Range("I25").Value = 0.9781962888247904
$ws.Range("J25").Value = 0.14526294131079
$ws.Range("N25").Value = 1.28741655399243
